$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values (prices, volume % changes, and row10/row11 coin swap)
$ws.Range('D2').Value = '40.020.78'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '2.210.23'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '294.45'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '87.32'
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('D7').Value = '0.513'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.470'
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '51.38'
$ws.Range('E10').Value = '  +7.32%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').Value = '30.66'
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '6.38'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = '2.555.77'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '13.80'
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = '2.196.55'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').Value = '39.965.22'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').Value = '65.45'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '235.11'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('D28').Value = '23.08'
$ws.Range('E28').Value = '  +2.31%  '
$ws.Range('D29').Value = '2.09'
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('D30').Value = '9.29'
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('D31').Value = '159.46'
$ws.Range('E31').Value = '  +2.33%  '
$ws.Range('D32').Value = '31.59'
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  +5.84%  '
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('D41').Value = '15.50'
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('D42').Value = '2.069.08'
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('D44').Value = '19.24'
$ws.Range('E44').Value = '  +10.46%  '
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').Value = '9.87'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('E47').Value = '  +3.21%  '
$ws.Range('D48').Value = '1.94'
$ws.Range('E48').Value = '  -9.58%  '
$ws.Range('D49').Value = '2.429.00'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '1.11'
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('E51').Value = '  +0.23%  '
